$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "base" (non-percentage) headers to their "_eig" counterparts.
# full -> full_eig, pruned -> pruned_eig, full_downsampled -> full_downsampled_eig,
# pruned_downsampled -> pruned_downsampled_eig
$ws.Range("A1").Value = "full_eig"
$ws.Range("C1").Value = "pruned_eig"
$ws.Range("E1").Value = "full_downsampled_eig"
$ws.Range("G1").Value = "pruned_downsampled_eig"

# New custom width for column G (pruned_downsampled_eig column header now needs more room).
# Note: Excel's ColumnWidth (character units) differs from the stored XML width; 11.71875
# is the COM value that round-trips to a stored width of 12.5, matching the target file.
$ws.Columns.Item(7).ColumnWidth = 11.71875

# Update the selected cell to I5 (reflects where the user was last working)
$ws.Range("I5").Select()
